# LxMLS Sponsorship Guidelines — 2023 -> 2025 edition update
$d = $word.ActiveDocument

# 1. Title: "LxMLS 2023 Sponsor Guide" -> "LxMLS 2025 Sponsor Guide"
$d.Content.Find.Execute("LxMLS 2023 Sponsor Guide", $true, $false, $false, $false, $false,
                         $true, 1, $false, "LxMLS 2025 Sponsor Guide", 2) | Out-Null

# 2. "more than 175 excellent students" -> "more than 200 excellent students"
$d.Content.Find.Execute("more than 175 excellent", $true, $false, $false, $false, $false,
                         $true, 1, $false, "more than 200 excellent", 2) | Out-Null

# 3. Rework the "In 2022 ... physical format" sentence opening (keeps the same
#    non-bold rPr formatting used throughout this run of text), folding in the
#    "Last year, through our sponsors ... 15 students" clause that used to
#    trail after the URL.
$old3 = "In 2022, the school was again organized in a physical format "
$new3 = "Since 2022, the school is being again organized in a physical format. Last year, through our sponsors, we fully or partially supported 15 students "
$d.Content.Find.Execute($old3, $true, $false, $false, $false, $false,
                         $true, 1, $false, $new3, 2) | Out-Null

# 4. Bump the school-website URL year: http://lxmls.it.pt/2022/ -> .../2024/
#    (kept as its own Find so the italic run formatting around the URL is preserved)
$d.Content.Find.Execute("http://lxmls.it.pt/2022/", $true, $false, $false, $false, $false,
                         $true, 1, $false, "http://lxmls.it.pt/2024/", 2) | Out-Null

# 5. Drop the now-redundant trailing clause (the "25 students" count moved
#    earlier in the sentence in step 3), leaving just the sentence terminator.
$old5 = " -- and, through our sponsors, we fully or partially supported 25 students. "
$new5 = ". "
$d.Content.Find.Execute($old5, $true, $false, $false, $false, $false,
                         $true, 1, $false, $new5, 2) | Out-Null

Write-Output "edits applied"
